$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 4
    7  = 7
    8  = 1
    9  = 4
    10 = 4
    11 = 2
    12 = 3
    13 = 3
    14 = 3
    15 = 3
    16 = 4
    17 = 2
    18 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
